$d = $word.ActiveDocument

# Locate the "SectionHeadnote" paragraph that currently reads
# "What is a corporation?" -- it is the first paragraph with that style,
# and it will become the sole survivor that carries the merged text.
$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("What is a corporation?")) {
        $startPara = $p
        break
    }
}

# Find the end of the final "SectionHeadnote" paragraph of Section Two,
# which currently reads "This is the second chapter of the casebook.".
$endPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("This is the second chapter of the casebook.")) {
        $endPara = $p
    }
}

# Delete every paragraph strictly after the "What is a corporation?"
# paragraph through (and including) the "This is the second chapter..."
# paragraph -- this removes all of the Resource/Case/Section
# sub-structure in between, merging it away.
$deleteRange = $d.Range($startPara.Range.End, $endPara.Range.End)
$deleteRange.Delete()

# Replace the remaining paragraph's own text with the concatenation of
# all of the text that used to live in the deleted paragraphs (the
# "This is an annotatable resource in the casebook." headnote text is
# not carried forward). A literal newline (not a paragraph mark) is
# embedded before "2Section Two" -- use `n so only one paragraph mark
# remains.
$newText = "1.1Case of the District Number 1This is the body of case 1.1.2Case of the District Number 2highlighted: content to highlight; elided: content to elide; replaced: content to replace; commented: content to comment; highlighted2: second highlight content;`n2Section Two"

$startPara.Range.Text = $newText
